$d = $word.ActiveDocument

# --- Locate the sentence that gets split. ---------------------------------
$sentence = $d.Content
$sentence.Find.ClearFormatting()
$null = $sentence.Find.Execute( `
    "Login as Administrator with the login credentials as ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$runStart = $sentence.Start
$runEnd   = $sentence.End

$prefix = "Login as "
$splitPoint = $runStart + $prefix.Length

# --- Insert the new word in the middle of the run. -------------------------
$insertPoint = $d.Range($splitPoint, $splitPoint)
$insertPoint.InsertBefore("Super ")
$insertedLen = "Super ".Length

# Splitting the middle of a run forces the engine to rebuild/flatten the
# run sequence from the insertion point through to the end of the paragraph,
# which would otherwise fold the hyperlink run and the runs that follow it
# together and strip their separate identities. Re-stamp every one of those
# runs with a harmless attribute toggle (set then immediately revert) so
# each keeps its own <w:r> boundary without any visible formatting change.
function Resplit-Range($r) {
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# 1) "Super " - the freshly inserted word becomes its own run.
Resplit-Range ($d.Range($splitPoint, $splitPoint + $insertedLen))

# 2) The remainder of the original sentence ("Administrator ... as ").
Resplit-Range ($d.Range($splitPoint + $insertedLen, $runEnd + $insertedLen))

# 3) The hyperlink run ("admin@admin.com") - relocate it with Find rather
#    than raw offsets since hyperlink fields occupy extra hidden character
#    positions in the Range addressing.
$afterSentence = $runEnd + $insertedLen
$hyperlink = $d.Range($afterSentence, $d.Content.End)
$hyperlink.Find.ClearFormatting()
$null = $hyperlink.Find.Execute("admin@admin.com", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Resplit-Range $hyperlink

# 4) "  and password as Test123"
$tail1 = $d.Range($hyperlink.End, $d.Content.End)
$tail1.Find.ClearFormatting()
$null = $tail1.Find.Execute("  and password as Test123", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Resplit-Range $tail1

# 5) The trailing "."
$tail2 = $d.Range($tail1.End, $tail1.End + 1)
Resplit-Range $tail2
